$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws4 = $wb.Worksheets.Item(4)

$ws1.Range("F2").Value = 2834
$ws1.Range("F3").Value = 1583
$ws1.Range("F5").Value = 557
$ws1.Range("F6").Value = 9587
$ws1.Range("F13").Value = 693
$ws1.Range("F14").Value = 693
$ws1.Range("F16").Value = 1201
$ws1.Range("F18").Value = 2982
$ws1.Range("F19").Value = 2254
$ws1.Range("F21").Value = 1955
$ws1.Range("F25").Value = 1563
$ws1.Range("F26").Value = 303
$ws1.Range("F27").Value = 19
$ws1.Range("F28").Value = 183
$ws1.Range("F30").Value = 33
$ws1.Range("F31").Value = 346
$ws1.Range("F34").Value = 517
$ws1.Range("F35").Value = 27
$ws1.Range("F36").Value = 133
$ws1.Range("F37").Value = 1541
$ws1.Range("F38").Value = 157
$ws1.Range("F39").Value = 1500
$ws1.Range("F40").Value = 37
$ws1.Range("F41").Value = 347
$ws1.Range("F42").Value = 27
$ws1.Range("F43").Value = 371
$ws1.Range("F44").Value = 762
$ws1.Range("F46").Value = 317
$ws2.Range("G2").Value = "不可售"
$ws4.Range("F2").Value = 2834
$ws4.Range("F3").Value = 1583
$ws4.Range("G4").Value = "不可售"
$ws4.Range("F5").Value = 557
$ws4.Range("F6").Value = 9587
$ws4.Range("F15").Value = 693
$ws4.Range("F16").Value = 693
$ws4.Range("F17").Value = 1201
$ws4.Range("F19").Value = 2982
$ws4.Range("F20").Value = 2254
$ws4.Range("F21").Value = 1955
$ws4.Range("F24").Value = 1563
$ws4.Range("F25").Value = 303
$ws4.Range("F26").Value = 19
$ws4.Range("F27").Value = 183
$ws4.Range("F29").Value = 33
$ws4.Range("F30").Value = 346
$ws4.Range("F33").Value = 517
$ws4.Range("F37").Value = 27
$ws4.Range("F38").Value = 133
$ws4.Range("F39").Value = 1541
$ws4.Range("F41").Value = 157
$ws4.Range("F42").Value = 1500
$ws4.Range("F43").Value = 37
$ws4.Range("F45").Value = 347
$ws4.Range("F46").Value = 27
$ws4.Range("F47").Value = 371
$ws4.Range("F48").Value = 762
$ws4.Range("F49").Value = 317
